# NatmiData Fgf8-Fgfrl1 LR-pairs workbook: "update scripts wuth new tpm"
#
# The sending cluster for this edge table changes from "MuSCs" to
# "Resolving-Mac" (column A, all data rows), and every TPM-derived
# expression / specificity metric is recomputed against the new TPM
# values. Target-cluster labels in column D are unaffected (MuSCs and
# Resolving-Mac still each appear once as a target, just on the other
# row than before).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: Sending cluster -> Resolving-Mac for every data row
$ws.Range("A2:A6").Value = "Resolving-Mac"

# Row 2 (Target cluster: ECs)
$ws.Range("G2").Value = 0.02530666666666667
$ws.Range("H2").Value = 0.07592
$ws.Range("M2").Value = 1.921622333333333
$ws.Range("N2").Value = 5.764867
$ws.Range("O2").Value = 0.1392241219313625
$ws.Range("P2").Value = 0.1392241219313625
$ws.Range("Q2").Value = 0.04862985584888889
$ws.Range("R2").Value = 0.43766870264
$ws.Range("S2").Value = 0.1392241219313625
$ws.Range("T2").Value = 0.1392241219313625

# Row 3 (Target cluster: FAPs)
$ws.Range("G3").Value = 0.02530666666666667
$ws.Range("H3").Value = 0.07592
$ws.Range("O3").Value = 0.7511588049189343
$ws.Range("P3").Value = 0.7511588049189343
$ws.Range("Q3").Value = 0.26237367416
$ws.Range("R3").Value = 2.36136306744
$ws.Range("S3").Value = 0.7511588049189343
$ws.Range("T3").Value = 0.7511588049189343

# Row 4 (Target cluster: Inflammatory-Mac)
$ws.Range("G4").Value = 0.02530666666666667
$ws.Range("H4").Value = 0.07592
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1441973333333333
$ws.Range("N4").Value = 0.432592
$ws.Range("O4").Value = 0.01044729069283506
$ws.Range("P4").Value = 0.01044729069283506
$ws.Range("Q4").Value = 0.003649153848888889
$ws.Range("R4").Value = 0.03284238464
$ws.Range("S4").Value = 0.01044729069283506
$ws.Range("T4").Value = 0.01044729069283506

# Row 5 (Target cluster: MuSCs)
$ws.Range("G5").Value = 0.02530666666666667
$ws.Range("H5").Value = 0.07592
$ws.Range("M5").Value = 1.182384
$ws.Range("N5").Value = 3.547152
$ws.Range("O5").Value = 0.08566531067535062
$ws.Range("P5").Value = 0.08566531067535062
$ws.Range("Q5").Value = 0.02992219776
$ws.Range("R5").Value = 0.26929977984
$ws.Range("S5").Value = 0.08566531067535062
$ws.Range("T5").Value = 0.08566531067535062

# Row 6 (Target cluster: Resolving-Mac)
$ws.Range("G6").Value = 0.02530666666666667
$ws.Range("H6").Value = 0.07592
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1863936666666667
$ws.Range("N6").Value = 0.559181
$ws.Range("O6").Value = 0.01350447178151746
$ws.Range("P6").Value = 0.01350447178151746
$ws.Range("Q6").Value = 0.004717002391111112
$ws.Range("R6").Value = 0.04245302152
$ws.Range("S6").Value = 0.01350447178151746
$ws.Range("T6").Value = 0.01350447178151746
